$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()
$newRow.HeadingFormat = $false

$leftQuote = [char]0x201C
$rightQuote = [char]0x201D

$newRow.Cells.Item(1).Range.Text = "January 4th 2022"
$newRow.Cells.Item(3).Range.Text = "- Minor change: added 6th option " + $leftQuote + "results" + $rightQuote + " to both flow chart and architectural design. (Architectural design and flow chart is now complete)"
